$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value2 = 1395.1904
$ws.Range("I98").Value2 = 1439.6842
$ws.Range("J98").Value2 = 972.5
$ws.Range("K98").Value2 = 1439.6842
$ws.Range("L98").Value2 = 972.5
$ws.Range("M98").Value2 = 58.31580000000008
$ws.Range("N98").Value2 = -3968.5

$ws.Range("H113").Value2 = 6043.5713
$ws.Range("I113").Value2 = 5751.25
$ws.Range("J113").Value2 = 6433.3335
$ws.Range("K113").Value2 = 5751.25
$ws.Range("L113").Value2 = 6433.3335
$ws.Range("M113").Value2 = -2497.25
$ws.Range("N113").Value2 = -12941.3335

$ws.Range("H122").Value2 = 1395.1904
$ws.Range("I122").Value2 = 1439.6842
$ws.Range("J122").Value2 = 972.5
$ws.Range("K122").Value2 = 4319.0526
$ws.Range("L122").Value2 = 2917.5
$ws.Range("M122").Value2 = -1869.0526
$ws.Range("N122").Value2 = -7817.5

$ws.Range("H131").Value2 = 1521.3636
$ws.Range("I131").Value2 = 970.55554
$ws.Range("K131").Value2 = 2911.66662
$ws.Range("M131").Value2 = 2128.33338

$ws.Range("H137").Value2 = 1642.4103
$ws.Range("I137").Value2 = 1568.32
$ws.Range("J137").Value2 = 1774.7142
$ws.Range("K137").Value2 = 4704.96
$ws.Range("L137").Value2 = 5324.142599999999
$ws.Range("M137").Value2 = -2154.96
$ws.Range("N137").Value2 = -10424.1426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 4078
$ws.Range("I2").Value2 = 4800
$ws.Range("J2").Value2 = 2995
$ws.Range("K2").Value2 = 4800
$ws.Range("L2").Value2 = 2995
$ws.Range("M2").Value2 = -4687
$ws.Range("N2").Value2 = -3221

$ws.Range("H32").Value2 = 4251.486
$ws.Range("I32").Value2 = 3490.7454
$ws.Range("K32").Value2 = 3490.7454
$ws.Range("M32").Value2 = -3203.7454

$ws.Range("H102").Value2 = 5293848
$ws.Range("I102").Value2 = 7409487
$ws.Range("J102").Value2 = 4750
$ws.Range("K102").Value2 = 7409487
$ws.Range("L102").Value2 = 4750
$ws.Range("M102").Value2 = -7407865
$ws.Range("N102").Value2 = -7994

$ws.Range("H116").Value2 = 4078
$ws.Range("I116").Value2 = 4800
$ws.Range("J116").Value2 = 2995
$ws.Range("K116").Value2 = 4800
$ws.Range("L116").Value2 = 2995
$ws.Range("M116").Value2 = -2506
$ws.Range("N116").Value2 = -7583

$ws.Range("H132").Value2 = 5298.6665
$ws.Range("I132").Value2 = 1541.3103
$ws.Range("J132").Value2 = 13680.462
$ws.Range("K132").Value2 = 4623.9309
$ws.Range("L132").Value2 = 41041.386
$ws.Range("M132").Value2 = -2093.9309
$ws.Range("N132").Value2 = -46101.386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 4078
$ws.Range("I3").Value2 = 4800
$ws.Range("J3").Value2 = 2995
$ws.Range("K3").Value2 = 4800
$ws.Range("L3").Value2 = 2995
$ws.Range("M3").Value2 = -4686
$ws.Range("N3").Value2 = -3223

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 18672.334
$ws.Range("I31").Value2 = 3007.4
$ws.Range("J31").Value2 = 38253.5
$ws.Range("K31").Value2 = 3007.4
$ws.Range("L31").Value2 = 38253.5
$ws.Range("M31").Value2 = -2712.4
$ws.Range("N31").Value2 = -38843.5

$ws.Range("H34").Value2 = 18672.334
$ws.Range("I34").Value2 = 3007.4
$ws.Range("J34").Value2 = 38253.5
$ws.Range("K34").Value2 = 3007.4
$ws.Range("L34").Value2 = 38253.5
$ws.Range("M34").Value2 = -2805.4
$ws.Range("N34").Value2 = -38657.5

$ws.Range("H86").Value2 = 2613.5334
$ws.Range("I86").Value2 = 2707.9167
$ws.Range("J86").Value2 = 2236
$ws.Range("K86").Value2 = 2707.9167
$ws.Range("L86").Value2 = 2236
$ws.Range("M86").Value2 = -1584.9167
$ws.Range("N86").Value2 = -4482

$ws.Range("H89").Value2 = 2613.5334
$ws.Range("I89").Value2 = 2707.9167
$ws.Range("J89").Value2 = 2236
$ws.Range("K89").Value2 = 13539.5835
$ws.Range("L89").Value2 = 11180
$ws.Range("M89").Value2 = -7923.583500000001
$ws.Range("N89").Value2 = -22412

$ws.Range("H94").Value2 = 2931.4167
$ws.Range("I94").Value2 = 2628.4546
$ws.Range("J94").Value2 = 3187.7693
$ws.Range("K94").Value2 = 2628.4546
$ws.Range("L94").Value2 = 3187.7693
$ws.Range("M94").Value2 = -2177.4546
$ws.Range("N94").Value2 = -4089.7693

$ws.Range("H132").Value2 = 3018.0667
$ws.Range("I132").Value2 = 2528.7778
$ws.Range("J132").Value2 = 3752
$ws.Range("K132").Value2 = 7586.3334
$ws.Range("L132").Value2 = 11256
$ws.Range("M132").Value2 = -5056.3334
$ws.Range("N132").Value2 = -16316

$ws.Range("H134").Value2 = 1544.875
$ws.Range("I134").Value2 = 1498.85
$ws.Range("J134").Value2 = 1775
$ws.Range("K134").Value2 = 4496.549999999999
$ws.Range("L134").Value2 = 5325
$ws.Range("M134").Value2 = -1961.549999999999
$ws.Range("N134").Value2 = -10395

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value2 = 22500346
$ws.Range("I12").Value2 = 16666947
$ws.Range("J12").Value2 = 31250444
$ws.Range("K12").Value2 = 50000841
$ws.Range("L12").Value2 = 93751332
$ws.Range("M12").Value2 = -50000668
$ws.Range("N12").Value2 = -93751678

$ws.Range("H92").Value2 = 699.2857
$ws.Range("I92").Value2 = 599
$ws.Range("J92").Value2 = 950
$ws.Range("K92").Value2 = 1797
$ws.Range("L92").Value2 = 2850
$ws.Range("M92").Value2 = -549
$ws.Range("N92").Value2 = -5346

$ws.Range("H132").Value2 = 2010.8195
$ws.Range("J132").Value2 = 2200.182
$ws.Range("L132").Value2 = 19801.638
$ws.Range("N132").Value2 = -24861.638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 5420.265
$ws.Range("J70").Value2 = 5453.643
$ws.Range("L70").Value2 = 5453.643
$ws.Range("N70").Value2 = -5993.643

$ws.Range("H73").Value2 = 5420.265
$ws.Range("J73").Value2 = 5453.643
$ws.Range("L73").Value2 = 5453.643
$ws.Range("N73").Value2 = -7325.643

$ws.Range("H102").Value2 = 3651.9412
$ws.Range("I102").Value2 = 1827.1818
$ws.Range("J102").Value2 = 6997.3335
$ws.Range("K102").Value2 = 1827.1818
$ws.Range("L102").Value2 = 6997.3335
$ws.Range("M102").Value2 = -205.1818000000001
$ws.Range("N102").Value2 = -10241.3335

$ws.Range("H132").Value2 = 3701.3438
$ws.Range("I132").Value2 = 4713.1816
$ws.Range("J132").Value2 = 3171.3333
$ws.Range("K132").Value2 = 14139.5448
$ws.Range("L132").Value2 = 9513.999899999999
$ws.Range("M132").Value2 = -11609.5448
$ws.Range("N132").Value2 = -14573.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value2 = 977.05884
$ws.Range("I16").Value2 = 593.0769
$ws.Range("J16").Value2 = 2225
$ws.Range("K16").Value2 = 593.0769
$ws.Range("L16").Value2 = 2225
$ws.Range("M16").Value2 = -423.0769
$ws.Range("N16").Value2 = -2565

$ws.Range("H122").Value2 = 3704343.5
$ws.Range("I122").Value2 = 4204897.5
$ws.Range("J122").Value2 = 2002460
$ws.Range("K122").Value2 = 12614692.5
$ws.Range("L122").Value2 = 6007380
$ws.Range("M122").Value2 = -12612242.5
$ws.Range("N122").Value2 = -6012280

$ws.Range("H136").Value2 = 11621.174
$ws.Range("I136").Value2 = 8639.333000000001
$ws.Range("K136").Value2 = 25917.999
$ws.Range("M136").Value2 = -23367.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value2 = 14011
$ws.Range("I39").Value2 = 2044
$ws.Range("J39").Value2 = 18000
$ws.Range("K39").Value2 = 2044
$ws.Range("L39").Value2 = 18000
$ws.Range("M39").Value2 = -1631
$ws.Range("N39").Value2 = -18826

$ws.Range("H113").Value2 = 709.175
$ws.Range("I113").Value2 = 664.28125
$ws.Range("K113").Value2 = 1992.84375
$ws.Range("M113").Value2 = 177.15625

$ws.Range("H122").Value2 = 1996.4138
$ws.Range("I122").Value2 = 1235.48
$ws.Range("J122").Value2 = 6752.25
$ws.Range("K122").Value2 = 3706.44
$ws.Range("L122").Value2 = 20256.75
$ws.Range("M122").Value2 = -1256.44
$ws.Range("N122").Value2 = -25156.75

$ws.Range("H123").Value2 = 29959.143
$ws.Range("J123").Value2 = 29959.143
$ws.Range("L123").Value2 = 29959.143
$ws.Range("N123").Value2 = -39759.143
